$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Populate cell values / hyperlinks in the exact order needed
#    so that the shared-strings table comes out in the same
#    order as the target workbook.
# ---------------------------------------------------------------

# Row 2: header
$ws.Range("A2").Value = 'Description'
$ws.Range("B2").Value = 'Example URL'

# Apply the Arial/9pt/wrap formatting to A2 now (before any
# hyperlinks are created) so this font/style is registered before
# the Hyperlink style, matching the target style-table ordering.
$ws.Range("A2").WrapText = $true
$ws.Range("A2").Font.Name = "Arial"
$ws.Range("A2").Font.Size = 9

# Rows 3-5: description then hyperlink, in order
$ws.Range("A3").Value = '3.2uF surface mount capacitors 0805'
$ws.Hyperlinks.Add($ws.Range("B3"), 'http://china.rs-online.com/web/p/ceramic-multilayer-capacitors/6911161/') | Out-Null
$ws.Range("A4").Value = '10k surface mount resistors 0805'
$ws.Hyperlinks.Add($ws.Range("B4"), 'http://china.rs-online.com/web/p/surface-mount-fixed-resistors/6789667/') | Out-Null
$ws.Range("A5").Value = '5k surface mount resistors 0805'
$ws.Hyperlinks.Add($ws.Range("B5"), 'http://china.rs-online.com/web/p/surface-mount-fixed-resistors/6791569/') | Out-Null

# Rows 6-9: descriptions entered first ...
$ws.Range("A6").Value = 'LM324G opamps'
$ws.Range("A7").Value = 'FMMT449 (NPN transistors)'
$ws.Range("A8").Value = 'FSB749 (PNP transistors)'
$ws.Range("A9").Value = 'DAC (AD5338RBRUZ)'

# ... then their hyperlinks added afterwards
$ws.Hyperlinks.Add($ws.Range("B6"), 'http://china.rs-online.com/web/p/operational-amplifiers/0858405/') | Out-Null
$ws.Hyperlinks.Add($ws.Range("B7"), 'http://china.rs-online.com/web/p/bipolar-transistors/6697681/') | Out-Null
$url8 = 'http://china.rs-online.com/web/p/bipolar-transistors/8076033/?searchTerm=FSB749&relevancy-data=636F3D3226696E3D4931384E4B6E6F776E41734D504E266C753D7A68266D6D3D6D61746368616C6C7061727469616C26706D3D5E5B5C772D5C2E2F252C5D2B2426706F3D313326736E3D592673743D4B4559574F52445F53494E474C455F414C5048415F4E554D455249432677633D424F5448267573743D465342373439267374613D46534237343926'
$ws.Hyperlinks.Add($ws.Range("B8"), $url8, "", "", $url8) | Out-Null
$url9 = 'http://china.rs-online.com/web/p/general-purpose-dacs/8209173/?searchTerm=AD5338RBRUZ&relevancy-data=636F3D3226696E3D4931384E4B6E6F776E41734D504E266C753D7A68266D6D3D6D61746368616C6C7061727469616C26706D3D5E5B5C772D5C2E2F252C5D2B2426706F3D313326736E3D592673743D4B4559574F52445F53494E474C455F414C5048415F4E554D455249432677633D424F5448267573743D414435333338524252555A267374613D414435333338524252555A26'
$ws.Hyperlinks.Add($ws.Range("B9"), $url9, "", "", $url9) | Out-Null

# Rows 10-11: description then hyperlink, in order
$ws.Range("A10").Value = '1k surface mount resistors 0805'
$ws.Hyperlinks.Add($ws.Range("B10"), 'http://china.rs-online.com/web/p/surface-mount-fixed-resistors/8145889/') | Out-Null
$ws.Range("A11").Value = '30k surface mount resistors 0805'
$ws.Hyperlinks.Add($ws.Range("B11"), 'http://china.rs-online.com/web/p/surface-mount-fixed-resistors/6792039/') | Out-Null

# C2 last
$ws.Range("C2").Value = 'Number'

# ---------------------------------------------------------------
# 2) Propagate A2 formatting (Arial 9pt, wrapped) to the rest of
#    column A and to B2 via copy/paste-special so we do not
#    regenerate duplicate font/style entries per cell.
# ---------------------------------------------------------------
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A3:A11").PasteSpecial(-4122) | Out-Null
$ws.Range("A2").Copy() | Out-Null
$ws.Range("B2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------
# 3) Column widths / row heights
# ---------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 21.592447916666668
$ws.Columns("B").ColumnWidth = 96.16666666666667
$ws.Columns("C").ColumnWidth = 34.022135416666664

$ws.Rows(3).RowHeight = 30
$ws.Rows(4).RowHeight = 30
$ws.Rows(5).RowHeight = 24.75
$ws.Rows(7).RowHeight = 24.75
$ws.Rows(8).RowHeight = 60
$ws.Rows(9).RowHeight = 60
$ws.Rows(10).RowHeight = 24.75
$ws.Rows(11).RowHeight = 24.75

# ---------------------------------------------------------------
# 4) Page setup + final selection
# ---------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("C4").Select() | Out-Null

Write-Output "done"
